$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "attendance"

# Update header row
$ws.Range("C1").Value = "19_05_18"
$ws.Range("D1").Value = "24_05_18"

# Update row 3
$ws.Range("A3").Value = "B17085"
$ws.Range("B3").Value = "yoy"

# Update row 4
$ws.Range("A4").Value = "B17110"
$ws.Range("B4").Value = "Varun Singh"

# Update row 5
$ws.Range("A5").Value = "b17056"
$ws.Range("B5").Value = "mea"
$ws.Range("D5").Value = 1

# Add new row 6
$ws.Range("A6").Value = "b17099"
$ws.Range("B6").Value = "me"
